$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks numeric but must stay literal text
# (decimal numbers with a trailing zero that Excel would otherwise normalize away).
# Mark them as Text-formatted before assigning so the literal string round-trips.
$textCells = @("D13", "D19", "D32", "D34", "D36", "D37", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.453.93"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.823.95"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "316.12"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "0.5171"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("D8").Value = "0.3854"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("D9").Value = "0.08272"
$ws.Range("E9").Value = "  +8.26%  "
$ws.Range("D10").Value = "1.122"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("D11").Value = "41.86"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "6.381"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").Value = "21.10"
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("D14").Value = "1.003"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "7.475"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").Value = "1.826.48"
$ws.Range("D17").Value = "93.97"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("E18").Value = "  +3.28%  "
$ws.Range("D19").Value = "0.06630"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").Value = "17.78"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "6.046"
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").Value = "28.490.25"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "11.46"
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").Value = "2.248"
$ws.Range("D26").Value = "21.09"
$ws.Range("E26").Value = "  +2.19%  "
$ws.Range("D27").Value = "159.37"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").Value = "2.034.99"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").Value = "0.1106"
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("D32").Value = "1.090"
$ws.Range("E32").Value = "  -3.28%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.723"
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.07540"
$ws.Range("E34").Value = "  +7.43%  "
$ws.Range("D35").Value = "3.684"
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("D36").Value = "0.2220"
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").Value = "0.02360"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("D38").Value = "12.06"
$ws.Range("E38").Value = "  +7.34%  "
$ws.Range("D39").Value = "5.242"
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("D40").Value = "8.755"
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("D41").Value = "0.6388"
$ws.Range("E41").Value = "  +1.78%  "
$ws.Range("D42").Value = "1.188"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").Value = "1.394"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").Value = "13.62"
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("D45").Value = "0.6191"
$ws.Range("E45").Value = "  +4.81%  "
$ws.Range("D46").Value = "3.795"
$ws.Range("E46").Value = "  +2.15%  "
$ws.Range("D47").Value = "127.49"
$ws.Range("E47").Value = "  +2.39%  "
$ws.Range("D48").Value = "2.004"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("D49").Value = "1.203"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "0.06960"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").Value = "1.080"
$ws.Range("E51").Value = "  +1.10%  "

# Restore the default (General) number format on those cells now that the
# literal text has been committed, so no stray cell style lingers behind.
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "general"
}
